# Weekly update for 'Fruta, Terminal La Palmera de La Serena - Pina' sheet.
# A new week's price block (4 sub-rows: Especial/Primera/Segunda/Tercera) is
# inserted at the top of the date-ordered data table (rows 473-476). Every
# existing block of rows shifts down by one block (4 rows), and the oldest
# existing block (previously rows 555-558) is appended as brand-new rows
# 559-562 at the bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns that are constant across the whole data block; needed to
#     populate the brand-new rows 559-562 in full. ---
$constA = 8
$constB = 'Terminal La Palmera de La Serena'
$constC = 'Coquimbo'
$constE = 4
$constF = 'Fruta'
$constG = 100108
$constH = 'Tropicales y subtropicales'
$constI = 100108005
$constJ = 'Piña'
$constK = 'Caramelo'
$constR = 'Ecuador'

# --- Per-row data: Date(D), Grade(L), UnitsPerBox(M), Min(N), Max(O), Avg(P),
#     PackagingText(Q), Avg$/unit(S), UnitsNumber(T). Row numbers are the FINAL
#     (post-edit) row numbers 473-562.
$rows = @(
  @{ R=473; D=44504; L='Especial'; M=216; N=20000; O=21000; P=20500; Q='$/caja 10 unidades'; S=2050; T=10 }
  @{ R=474; D=44504; L='Primera'; M=216; N=20000; O=21000; P=20500; Q='$/caja 12 unidades'; S=1708; T=12 }
  @{ R=475; D=44504; L='Segunda'; M=216; N=20000; O=21000; P=20500; Q='$/caja 14 unidades'; S=1464; T=14 }
  @{ R=476; D=44504; L='Tercera'; M=216; N=20000; O=21000; P=20500; Q='$/caja 16 unidades'; S=1281; T=16 }
  @{ R=477; D=44246; L='Especial'; M=216; N=17000; O=17500; P=17250; Q='$/caja 10 unidades'; S=1725; T=10 }
  @{ R=478; D=44246; L='Primera'; M=208; N=17000; O=17500; P=17260; Q='$/caja 12 unidades'; S=1438; T=12 }
  @{ R=479; D=44246; L='Segunda'; M=216; N=17000; O=17500; P=17250; Q='$/caja 14 unidades'; S=1232; T=14 }
  @{ R=480; D=44487; L='Especial'; M=216; N=22500; O=23000; P=22750; Q='$/caja 10 unidades'; S=2275; T=10 }
  @{ R=481; D=44487; L='Primera'; M=216; N=22500; O=23000; P=22750; Q='$/caja 12 unidades'; S=1896; T=12 }
  @{ R=482; D=44487; L='Segunda'; M=216; N=22500; O=23000; P=22750; Q='$/caja 14 unidades'; S=1625; T=14 }
  @{ R=483; D=44487; L='Tercera'; M=216; N=22500; O=23000; P=22750; Q='$/caja 16 unidades'; S=1422; T=16 }
  @{ R=484; D=44425; L='Especial'; M=432; N=18500; O=19000; P=18750; Q='$/caja 10 unidades'; S=1875; T=10 }
  @{ R=485; D=44425; L='Primera'; M=432; N=18500; O=19000; P=18750; Q='$/caja 12 unidades'; S=1562; T=12 }
  @{ R=486; D=44425; L='Segunda'; M=432; N=18500; O=19000; P=18750; Q='$/caja 14 unidades'; S=1339; T=14 }
  @{ R=487; D=44343; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=488; D=44343; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=489; D=44343; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=490; D=44343; L='Tercera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 16 unidades'; S=1047; T=16 }
  @{ R=491; D=44449; L='Especial'; M=216; N=22000; O=23000; P=22500; Q='$/caja 10 unidades'; S=2250; T=10 }
  @{ R=492; D=44449; L='Primera'; M=216; N=22000; O=23000; P=22500; Q='$/caja 12 unidades'; S=1875; T=12 }
  @{ R=493; D=44449; L='Segunda'; M=216; N=22000; O=23000; P=22500; Q='$/caja 14 unidades'; S=1607; T=14 }
  @{ R=494; D=44449; L='Tercera'; M=216; N=22000; O=23000; P=22500; Q='$/caja 16 unidades'; S=1406; T=16 }
  @{ R=495; D=44168; L='Especial'; M=216; N=23500; O=24000; P=23750; Q='$/caja 10 unidades'; S=2375; T=10 }
  @{ R=496; D=44168; L='Primera'; M=216; N=23500; O=24000; P=23750; Q='$/caja 12 unidades'; S=1979; T=12 }
  @{ R=497; D=44168; L='Segunda'; M=216; N=23500; O=24000; P=23750; Q='$/caja 14 unidades'; S=1696; T=14 }
  @{ R=498; D=44175; L='Especial'; M=216; N=23000; O=23500; P=23250; Q='$/caja 10 unidades'; S=2325; T=10 }
  @{ R=499; D=44175; L='Primera'; M=216; N=23000; O=23500; P=23250; Q='$/caja 12 unidades'; S=1938; T=12 }
  @{ R=500; D=44175; L='Segunda'; M=216; N=23000; O=23500; P=23250; Q='$/caja 14 unidades'; S=1661; T=14 }
  @{ R=501; D=44392; L='Especial'; M=216; N=18500; O=19000; P=18750; Q='$/caja 10 unidades'; S=1875; T=10 }
  @{ R=502; D=44392; L='Primera'; M=216; N=18500; O=19000; P=18750; Q='$/caja 12 unidades'; S=1562; T=12 }
  @{ R=503; D=44392; L='Segunda'; M=216; N=18500; O=19000; P=18750; Q='$/caja 14 unidades'; S=1339; T=14 }
  @{ R=504; D=44473; L='Especial'; M=216; N=20000; O=21000; P=20500; Q='$/caja 10 unidades'; S=2050; T=10 }
  @{ R=505; D=44473; L='Primera'; M=216; N=20000; O=21000; P=20500; Q='$/caja 12 unidades'; S=1708; T=12 }
  @{ R=506; D=44473; L='Segunda'; M=216; N=20000; O=21000; P=20500; Q='$/caja 14 unidades'; S=1464; T=14 }
  @{ R=507; D=44400; L='Especial'; M=216; N=19500; O=20000; P=19750; Q='$/caja 10 unidades'; S=1975; T=10 }
  @{ R=508; D=44400; L='Primera'; M=208; N=19500; O=20000; P=19760; Q='$/caja 12 unidades'; S=1647; T=12 }
  @{ R=509; D=44400; L='Segunda'; M=216; N=19500; O=20000; P=19750; Q='$/caja 14 unidades'; S=1411; T=14 }
  @{ R=510; D=44484; L='Especial'; M=216; N=23000; O=24000; P=23500; Q='$/caja 10 unidades'; S=2350; T=10 }
  @{ R=511; D=44484; L='Primera'; M=216; N=23000; O=24000; P=23500; Q='$/caja 12 unidades'; S=1958; T=12 }
  @{ R=512; D=44484; L='Segunda'; M=216; N=23000; O=24000; P=23500; Q='$/caja 14 unidades'; S=1679; T=14 }
  @{ R=513; D=44484; L='Tercera'; M=216; N=23000; O=24000; P=23500; Q='$/caja 16 unidades'; S=1469; T=16 }
  @{ R=514; D=44494; L='Especial'; M=216; N=21000; O=22000; P=21500; Q='$/caja 10 unidades'; S=2150; T=10 }
  @{ R=515; D=44494; L='Primera'; M=216; N=21000; O=22000; P=21500; Q='$/caja 12 unidades'; S=1792; T=12 }
  @{ R=516; D=44494; L='Segunda'; M=216; N=21000; O=22000; P=21500; Q='$/caja 14 unidades'; S=1536; T=14 }
  @{ R=517; D=44494; L='Tercera'; M=216; N=21000; O=22000; P=21500; Q='$/caja 16 unidades'; S=1344; T=16 }
  @{ R=518; D=44445; L='Especial'; M=216; N=20000; O=21000; P=20500; Q='$/caja 10 unidades'; S=2050; T=10 }
  @{ R=519; D=44445; L='Primera'; M=216; N=20000; O=21000; P=20500; Q='$/caja 12 unidades'; S=1708; T=12 }
  @{ R=520; D=44445; L='Segunda'; M=216; N=20000; O=21000; P=20500; Q='$/caja 14 unidades'; S=1464; T=14 }
  @{ R=521; D=44301; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=522; D=44301; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=523; D=44301; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=524; D=44330; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=525; D=44330; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=526; D=44330; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=527; D=44330; L='Tercera'; M=208; N=16500; O=17000; P=16760; Q='$/caja 16 unidades'; S=1048; T=16 }
  @{ R=528; D=44270; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=529; D=44270; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=530; D=44270; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=531; D=44295; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=532; D=44295; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=533; D=44295; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=534; D=44295; L='Tercera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 16 unidades'; S=1047; T=16 }
  @{ R=535; D=44217; L='Especial'; M=216; N=17500; O=18000; P=17750; Q='$/caja 10 unidades'; S=1775; T=10 }
  @{ R=536; D=44217; L='Primera'; M=216; N=17500; O=18000; P=17750; Q='$/caja 12 unidades'; S=1479; T=12 }
  @{ R=537; D=44217; L='Segunda'; M=216; N=17500; O=18000; P=17750; Q='$/caja 14 unidades'; S=1268; T=14 }
  @{ R=538; D=44421; L='Especial'; M=216; N=19000; O=19500; P=19250; Q='$/caja 10 unidades'; S=1925; T=10 }
  @{ R=539; D=44421; L='Primera'; M=216; N=19000; O=19500; P=19250; Q='$/caja 12 unidades'; S=1604; T=12 }
  @{ R=540; D=44421; L='Segunda'; M=216; N=19000; O=19500; P=19250; Q='$/caja 14 unidades'; S=1375; T=14 }
  @{ R=541; D=44421; L='Tercera'; M=216; N=19000; O=19500; P=19250; Q='$/caja 16 unidades'; S=1203; T=16 }
  @{ R=542; D=44383; L='Especial'; M=216; N=17500; O=18000; P=17750; Q='$/caja 10 unidades'; S=1775; T=10 }
  @{ R=543; D=44383; L='Primera'; M=216; N=17500; O=18000; P=17750; Q='$/caja 12 unidades'; S=1479; T=12 }
  @{ R=544; D=44383; L='Segunda'; M=216; N=17500; O=18000; P=17750; Q='$/caja 14 unidades'; S=1268; T=14 }
  @{ R=545; D=44273; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=546; D=44273; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=547; D=44273; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=548; D=44273; L='Tercera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 16 unidades'; S=1047; T=16 }
  @{ R=549; D=44433; L='Especial'; M=216; N=19500; O=20000; P=19750; Q='$/caja 10 unidades'; S=1975; T=10 }
  @{ R=550; D=44433; L='Primera'; M=216; N=19500; O=20000; P=19750; Q='$/caja 12 unidades'; S=1646; T=12 }
  @{ R=551; D=44433; L='Segunda'; M=216; N=19500; O=20000; P=19750; Q='$/caja 14 unidades'; S=1411; T=14 }
  @{ R=552; D=44433; L='Tercera'; M=216; N=19500; O=20000; P=19750; Q='$/caja 16 unidades'; S=1234; T=16 }
  @{ R=553; D=44302; L='Especial'; M=216; N=16500; O=17000; P=16750; Q='$/caja 10 unidades'; S=1675; T=10 }
  @{ R=554; D=44302; L='Primera'; M=216; N=16500; O=17000; P=16750; Q='$/caja 12 unidades'; S=1396; T=12 }
  @{ R=555; D=44302; L='Segunda'; M=216; N=16500; O=17000; P=16750; Q='$/caja 14 unidades'; S=1196; T=14 }
  @{ R=556; D=44179; L='Especial'; M=216; N=23000; O=23500; P=23250; Q='$/caja 10 unidades'; S=2325; T=10 }
  @{ R=557; D=44179; L='Primera'; M=216; N=23000; O=23500; P=23250; Q='$/caja 12 unidades'; S=1938; T=12 }
  @{ R=558; D=44179; L='Segunda'; M=216; N=23000; O=23500; P=23250; Q='$/caja 14 unidades'; S=1661; T=14 }
  @{ R=559; D=44491; L='Especial'; M=216; N=22000; O=22500; P=22250; Q='$/caja 10 unidades'; S=2225; T=10 }
  @{ R=560; D=44491; L='Primera'; M=216; N=22000; O=22500; P=22250; Q='$/caja 12 unidades'; S=1854; T=12 }
  @{ R=561; D=44491; L='Segunda'; M=216; N=22000; O=22500; P=22250; Q='$/caja 14 unidades'; S=1589; T=14 }
  @{ R=562; D=44491; L='Tercera'; M=216; N=22000; O=22500; P=22250; Q='$/caja 16 unidades'; S=1391; T=16 }
)

foreach ($row in $rows) {
  $r = $row.R
  if ($r -gt 558) {
    # Brand-new row: populate every column, including the ones that are
    # constant throughout the sheet.
    $ws.Cells.Item($r, 1).Value = $constA
    $ws.Cells.Item($r, 2).Value = $constB
    $ws.Cells.Item($r, 3).Value = $constC
    $ws.Cells.Item($r, 5).Value = $constE
    $ws.Cells.Item($r, 6).Value = $constF
    $ws.Cells.Item($r, 7).Value = $constG
    $ws.Cells.Item($r, 8).Value = $constH
    $ws.Cells.Item($r, 9).Value = $constI
    $ws.Cells.Item($r, 10).Value = $constJ
    $ws.Cells.Item($r, 11).Value = $constK
    $ws.Cells.Item($r, 18).Value = $constR
  }
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $ws.Cells.Item($r, 12).Value = $row.L
  $ws.Cells.Item($r, 13).Value = $row.M
  $ws.Cells.Item($r, 14).Value = $row.N
  $ws.Cells.Item($r, 15).Value = $row.O
  $ws.Cells.Item($r, 16).Value = $row.P
  $ws.Cells.Item($r, 17).Value = $row.Q
  $ws.Cells.Item($r, 19).Value = $row.S
  $ws.Cells.Item($r, 20).Value = $row.T
}

